$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: new customer entry
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "testcustomer@gmail.com"
$ws.Cells.Item(6, 3).Value = "Test"
$ws.Cells.Item(6, 4).Value = "Customer"
$ws.Cells.Item(6, 5).Value = "1st Street, Swansea, SA15BD"
$ws.Cells.Item(6, 6).Value = "Customer"

# Row 7: new customer entry
$ws.Cells.Item(7, 1).Value = 16
$ws.Cells.Item(7, 2).Value = "testsustomer1@gmail.com"
$ws.Cells.Item(7, 3).Value = "Test"
$ws.Cells.Item(7, 4).Value = "Test"
$ws.Cells.Item(7, 5).Value = "Test"
$ws.Cells.Item(7, 6).Value = "Customer"

# Widen column E to fit the new longer address text (best-fit for longest value)
$ws.Columns.Item(5).ColumnWidth = 25.75
